$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-unneeded rows 11-16 (old Post-Intervention sample data)
$ws.Range("A11:B16").EntireRow.Delete() | Out-Null

# Replace the sample texts in B2:B10 with the new content
$ws.Range("B2").Value = "I hope you are well. `nI have amended page 13 of the Terms of Business to include the sentence: “I confirm that I have authority to instruct solicitors on behalf of the company.” `nCould you please sign and return the updated Terms of Business, and also please sign and return the Letter of Engagement.  `nIf you have any questions please do not hesitate to contact me.  "
$ws.Range("B3").Value = "As instructed I have drafted the attached letters which will be posted <DATE>. `nHave a good weekend "
$ws.Range("B4").Value = "I hope you are well.  `nI have not received a response to our Letter of Claim from Ms J or Mr R.`nI can write to both opponents to advise that I have not received their response and that I look forward to receiving their response within 24 <DOMAIN> they fail to respond, you have the option to issue court proceedings against them which I can discuss with you.  `nWould you like me to send a letter to them? "
$ws.Range("B5").Value = "Just to confirm following our telephone conversation, I have a few quick questions to ask you in relation to the company structure.  `nThank you for agreeing to speak to me <DATE> at 12.15. `nHave a lovely evening. "
$ws.Range("B6").Value = "I have spoken with my costs team and we do not need you to sign a new <ORGANIZATION>. I need to make you aware however that the <ORGANIZATION1> included a reference to qualified one way cost shifting in the event that you were unsuccessful in the claim, this does not apply to the case and if you were to be unsuccessful, you would be liable to pay the Defendant's costs to be assessed if not agreed.`nAs to my fees, these are covered by the CFA with the agreed success fee. All work is carried out at my hourly rate of <MONEY> plus VAT. Our costs are presently <MONEY1> plus VAT. I anticipate the fees to trial could exceed <MONEY2> plus VAT. I will update you as to the costs at various points throughout the case however they are for information only, we do not require you to pay this."
$ws.Range("B7").Value = "Please find attached invoice for the caveat removal."
$ws.Range("B8").Value = "Out conveyancing team are dealing and will be liaising with the solicitors."
$ws.Range("B9").Value = "It doesn’t appear that’s what they have done, their emails suggest otherwise. `nMy real concern is this will all now fall apart. If P tried to engage with T and she has refused that is one thing but if they have been contacting her to say you are dealing, then she is right to say no. It’s difficult to say without sight of their paperwork.  `nShe can request your mums bank details to make payment but you need to know how much she has so would will need to speak with her to ensure all beneficiaries paid and she doesn’t keep the money which was the whole point in agreeing an independent solicitor. "
$ws.Range("B10").Value = "Sometimes technology is great, sometimes not. `nI can do <DATE> morning in our Walker House <ORGANIZATION1> 10-1 or <DATE1> <DATE2> or 1-3 in our <ORGANIZATION> office. `nPlease let me know if either of these work or I can ask my IT team to set up a 3 way call, I am unable to do it directly from my mobile. `nIn the meantime, I will continue my review so that when we speak, I will hopefully be able to answer some of your initial questions. "

# Re-fit rows so that embedded line breaks don't leave a stray custom row-height
# (the source file has no explicit row heights set).
$ws.Rows("1:10").AutoFit() | Out-Null
